# Auto-generated Excel COM-interop script
# Applies: (1) remove spaces around '-' in column E time ranges
#          (2) update column F (want-to-go counts) per diff

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("E2").Value = '2024.03.30 10:00-04.29 22:00'
$ws.Range("F2").Value = 256
$ws.Range("E3").Value = '2024.04.12 10:00-05.12 20:00'
$ws.Range("F3").Value = 878
$ws.Range("E4").Value = '2024.04.13 11:00-04.14 17:00'
$ws.Range("F4").Value = 552
$ws.Range("E5").Value = '2024.04.13 10:30-04.14 16:30'
$ws.Range("E6").Value = '2024.04.13 11:00-04.14 18:00'
$ws.Range("F6").Value = 1367
$ws.Range("E7").Value = '2024.04.13 10:00-04.14 18:00'
$ws.Range("F7").Value = 121
$ws.Range("E8").Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range("F8").Value = 815
$ws.Range("E9").Value = '2024.04.13 10:00-04.21 17:00'
$ws.Range("F9").Value = 1169
$ws.Range("E10").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E11").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("F11").Value = 3056
$ws.Range("E12").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E13").Value = '2024.04.20 13:50-04.20 18:00'
$ws.Range("E14").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E15").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("F15").Value = 623
$ws.Range("E16").Value = '2024.04.27 12:40-04.27 16:40'
$ws.Range("E17").Value = '2024.04.27 10:00-04.27 18:00'
$ws.Range("E18").Value = '2024.04.27 10:30-04.27 16:30'
$ws.Range("F18").Value = 613
$ws.Range("E19").Value = '2024.05.01 10:00-06.02 22:00'
$ws.Range("F19").Value = 1124
$ws.Range("E20").Value = '2024.05.01 10:00-06.02 22:00'
$ws.Range("F20").Value = 1124
$ws.Range("E21").Value = '2024.05.01 10:00-05.02 17:00'
$ws.Range("F21").Value = 161
$ws.Range("E22").Value = '2024.05.01 10:30-05.01 16:30'
$ws.Range("F22").Value = 538
$ws.Range("E23").Value = '2024.05.01 10:00-05.05 16:00'
$ws.Range("F23").Value = 191
$ws.Range("E24").Value = '2024.05.01 10:00-05.01 17:00'
$ws.Range("F24").Value = 538
$ws.Range("E25").Value = '2024.05.02 10:30-05.04 19:00'
$ws.Range("F25").Value = 251
$ws.Range("E26").Value = '2024.05.02 13:00-05.02 18:00'
$ws.Range("E27").Value = '2024.05.02 10:20-05.03 16:30'
$ws.Range("F27").Value = 607
$ws.Range("E28").Value = '2024.05.03 10:00-05.04 16:00'
$ws.Range("E29").Value = '2024.05.04 10:00-05.05 17:00'
$ws.Range("F29").Value = 842
$ws.Range("E30").Value = '2024.05.04 10:00-05.04 18:00'
$ws.Range("F30").Value = 86
$ws.Range("E31").Value = '2024.05.05 10:00-05.05 17:00'
$ws.Range("E32").Value = '2024.05.05 10:00-05.05 18:00'
$ws.Range("F32").Value = 53
$ws.Range("E33").Value = '2024.05.18 10:00-05.19 17:00'
$ws.Range("F33").Value = 1054
$ws.Range("E34").Value = '2024.05.18 10:00-05.19 17:00'
$ws.Range("F34").Value = 5065
$ws.Range("E35").Value = '2024.05.18 10:00-05.18 17:00'
$ws.Range("F35").Value = 525
$ws.Range("E36").Value = '2024.05.19 10:00-05.19 17:00'
$ws.Range("F36").Value = 259
$ws.Range("E37").Value = '2024.06.01 10:00-06.02 17:00'
$ws.Range("F37").Value = 138
$ws.Range("E38").Value = '2024.06.08 10:00-06.10 16:00'
$ws.Range("E39").Value = '2024.07.05 10:00-07.07 16:00'
$ws.Range("F39").Value = 6
$ws.Range("E40").Value = '2024.07.12 10:00-07.14 16:00'
$ws.Range("F40").Value = 6

$ws = $wb.Worksheets.Item(2)
$ws.Range("E2").Value = '2024.04.05 15:50-05.03 20:10'
$ws.Range("F2").Value = 18
$ws.Range("E3").Value = '2024.04.12 19:30-04.12 21:30'
$ws.Range("E4").Value = '2024.04.12 19:00-04.12 20:30'
$ws.Range("F4").Value = 367
$ws.Range("E5").Value = '2024.04.13 16:00-04.13 18:00'
$ws.Range("E6").Value = '2024.04.13 19:30-04.13 21:30'
$ws.Range("F6").Value = 412
$ws.Range("E7").Value = '2024.04.13 19:30-04.13 21:00'
$ws.Range("E8").Value = '2024.04.13 18:30-04.20 21:50'
$ws.Range("E9").Value = '2024.04.14 14:00-04.14 15:40'
$ws.Range("E10").Value = '2024.04.14 15:00-05.01 20:15'
$ws.Range("E11").Value = '2024.04.20 19:30-04.20 21:00'
$ws.Range("E12").Value = '2024.04.20 19:30-04.20 21:30'
$ws.Range("E13").Value = '2024.04.20 13:00-04.20 15:00'
$ws.Range("E14").Value = '2024.04.21 15:00-04.21 16:30'
$ws.Range("E15").Value = '2024.04.21 20:00-04.21 21:30'
$ws.Range("E17").Value = '2024.04.26 19:30-04.26 21:30'
$ws.Range("F17").Value = 64
$ws.Range("E18").Value = '2024.04.26 19:00-04.26 20:30'
$ws.Range("E19").Value = '2024.04.27 20:00-04.27 21:30'
$ws.Range("E20").Value = '2024.05.01 20:00-05.01 22:00'
$ws.Range("E21").Value = '2024.05.01 18:30-05.01 21:00'
$ws.Range("E22").Value = '2024.05.01 19:30-05.19 21:00'
$ws.Range("E23").Value = '2024.05.02 13:30-05.02 15:20'
$ws.Range("E24").Value = '2024.05.02 19:30-05.02 21:00'
$ws.Range("F24").Value = 51
$ws.Range("E25").Value = '2024.05.03 19:00-05.03 22:00'
$ws.Range("F25").Value = 388
$ws.Range("E26").Value = '2024.05.03 19:30-05.03 21:00'
$ws.Range("E27").Value = '2024.05.04 14:20-06.09 15:35'
$ws.Range("E28").Value = '2024.05.04 14:00-05.04 16:00'
$ws.Range("F28").Value = 692
$ws.Range("E29").Value = '2024.05.04 19:30-05.04 21:00'
$ws.Range("E30").Value = '2024.05.04 16:30-06.02 17:50'
$ws.Range("E31").Value = '2024.05.05 13:00-05.05 15:30'
$ws.Range("E32").Value = '2024.05.17 19:30-05.17 21:00'
$ws.Range("E33").Value = '2024.05.18 19:30-05.18 21:00'
$ws.Range("E34").Value = '2024.05.18 14:00-05.18 20:30'
$ws.Range("F34").Value = 62
$ws.Range("E35").Value = '2024.05.19 19:30-05.19 21:00'
$ws.Range("E36").Value = '2024.05.19 14:30-05.19 16:00'
$ws.Range("E37").Value = '2024.05.19 14:00-05.19 15:30'
$ws.Range("F37").Value = 444
$ws.Range("E38").Value = '2024.05.24 19:30-05.24 21:00'
$ws.Range("E39").Value = '2024.05.25 19:30-05.25 21:00'
$ws.Range("F39").Value = 13
$ws.Range("E40").Value = '2024.06.01 19:30-06.01 21:00'
$ws.Range("E41").Value = '2024.06.01 19:30-06.01 21:00'
$ws.Range("E42").Value = '2024.06.01 19:30-06.01 21:00'
$ws.Range("E43").Value = '2024.06.08 19:30-06.08 21:00'
$ws.Range("E44").Value = '2024.06.15 19:30-06.15 22:00'
$ws.Range("E45").Value = '2024.06.22 19:30-06.22 21:30'
$ws.Range("E46").Value = '2024.07.17 19:30-07.17 21:00'
$ws.Range("E47").Value = '2024.07.19 19:30-07.19 21:30'

$ws = $wb.Worksheets.Item(3)
$ws.Range("E2").Value = '2023.10.16 10:00-2024.10.15 21:00'
$ws.Range("E3").Value = '2023.10.25 10:00-2024.10.20 21:00'
$ws.Range("E4").Value = '2024.03.21 00:00-04.28 23:59'
$ws.Range("F4").Value = 644
$ws.Range("E5").Value = '2024.04.04 00:00-05.20 23:59'
$ws.Range("E6").Value = '2024.04.24 00:00-05.22 23:59'
$ws.Range("F6").Value = 419

$ws = $wb.Worksheets.Item(4)
$ws.Range("E2").Value = '2024.03.21 00:00-04.28 23:59'
$ws.Range("F2").Value = 644
$ws.Range("E3").Value = '2024.03.30 10:00-04.29 22:00'
$ws.Range("F3").Value = 256
$ws.Range("E4").Value = '2024.04.04 00:00-05.20 23:59'
$ws.Range("E5").Value = '2024.04.05 15:50-05.03 20:10'
$ws.Range("F5").Value = 18
$ws.Range("E6").Value = '2024.04.12 10:00-05.12 20:00'
$ws.Range("F6").Value = 878
$ws.Range("E7").Value = '2024.04.13 19:30-04.13 21:30'
$ws.Range("F7").Value = 412
$ws.Range("E8").Value = '2024.04.13 11:00-04.14 17:00'
$ws.Range("F8").Value = 552
$ws.Range("E9").Value = '2024.04.13 10:30-04.14 16:30'
$ws.Range("E10").Value = '2024.04.13 11:00-04.14 18:00'
$ws.Range("F10").Value = 1367
$ws.Range("E11").Value = '2024.04.13 10:00-04.14 18:00'
$ws.Range("F11").Value = 121
$ws.Range("E12").Value = '2024.04.13 10:00-04.13 17:00'
$ws.Range("F12").Value = 815
$ws.Range("E13").Value = '2024.04.13 10:00-04.21 17:00'
$ws.Range("F13").Value = 1169
$ws.Range("E14").Value = '2024.04.20 19:30-04.20 21:00'
$ws.Range("E15").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E16").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("F16").Value = 3056
$ws.Range("E17").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E18").Value = '2024.04.20 13:00-04.20 15:00'
$ws.Range("E19").Value = '2024.04.20 13:50-04.20 18:00'
$ws.Range("E20").Value = '2024.04.20 10:00-04.21 17:00'
$ws.Range("E21").Value = '2024.04.24 00:00-05.22 23:59'
$ws.Range("F21").Value = 419
$ws.Range("E22").Value = '2024.04.27 12:40-04.27 16:40'
$ws.Range("E23").Value = '2024.04.27 10:00-04.27 18:00'
$ws.Range("E24").Value = '2024.04.27 10:30-04.27 16:30'
$ws.Range("F24").Value = 613
$ws.Range("E25").Value = '2024.05.01 10:00-06.02 22:00'
$ws.Range("F25").Value = 1124
$ws.Range("E26").Value = '2024.05.01 10:00-06.02 22:00'
$ws.Range("F26").Value = 1124
$ws.Range("E27").Value = '2024.05.01 10:00-05.02 17:00'
$ws.Range("F27").Value = 161
$ws.Range("E28").Value = '2024.05.01 18:30-05.01 21:00'
$ws.Range("E29").Value = '2024.05.01 10:30-05.01 16:30'
$ws.Range("F29").Value = 538
$ws.Range("E30").Value = '2024.05.01 19:30-05.19 21:00'
$ws.Range("E31").Value = '2024.05.01 10:00-05.05 16:00'
$ws.Range("F31").Value = 191
$ws.Range("E32").Value = '2024.05.02 10:30-05.04 19:00'
$ws.Range("F32").Value = 251
$ws.Range("E33").Value = '2024.05.02 19:30-05.02 21:00'
$ws.Range("F33").Value = 51
$ws.Range("E34").Value = '2024.05.02 13:00-05.02 18:00'
$ws.Range("E35").Value = '2024.05.02 10:20-05.03 16:30'
$ws.Range("F35").Value = 607
$ws.Range("E36").Value = '2024.05.03 19:00-05.03 22:00'
$ws.Range("F36").Value = 388
$ws.Range("E37").Value = '2024.05.04 14:00-05.04 16:00'
$ws.Range("F37").Value = 692
$ws.Range("E38").Value = '2024.05.04 10:00-05.05 17:00'
$ws.Range("F38").Value = 842
$ws.Range("E39").Value = '2024.05.04 10:00-05.04 18:00'
$ws.Range("F39").Value = 86
$ws.Range("E40").Value = '2024.05.05 10:00-05.05 17:00'
$ws.Range("E41").Value = '2024.05.05 13:00-05.05 15:30'
$ws.Range("E42").Value = '2024.05.05 10:00-05.05 18:00'
$ws.Range("F42").Value = 53
$ws.Range("E43").Value = '2024.05.18 10:00-05.19 17:00'
$ws.Range("F43").Value = 1054
$ws.Range("E44").Value = '2024.05.18 10:00-05.19 17:00'
$ws.Range("F44").Value = 5065
$ws.Range("E45").Value = '2024.05.18 14:00-05.18 20:30'
$ws.Range("F45").Value = 62
$ws.Range("E46").Value = '2024.05.18 10:00-05.18 17:00'
$ws.Range("F46").Value = 525
$ws.Range("E47").Value = '2024.05.19 14:00-05.19 15:30'
$ws.Range("F47").Value = 444
$ws.Range("E48").Value = '2024.05.19 14:00-05.19 15:30'
$ws.Range("F48").Value = 444
$ws.Range("E49").Value = '2024.05.19 10:00-05.19 17:00'
$ws.Range("F49").Value = 259
$ws.Range("E50").Value = '2024.06.01 19:30-06.01 21:00'
$ws.Range("E51").Value = '2024.06.22 19:30-06.22 21:30'
$ws.Range("E52").Value = '2024.07.12 10:00-07.14 16:00'
$ws.Range("F52").Value = 6
